$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H26").Value = 0.6575
$ws.Range("I26").Value = 0.01883
$ws.Range("H27").Value = 0.18891
$ws.Range("I27").Value = 0.04049
$ws.Range("H28").Value = 0.63153
$ws.Range("I28").Value = 0.0162
$ws.Range("H29").Value = 0.08197
$ws.Range("I29").Value = 0.02841
$ws.Range("H30").Value = 0.65708
$ws.Range("I30").Value = 0.01873
$ws.Range("H31").Value = 0.18925
$ws.Range("I31").Value = 0.04184
$ws.Range("H32").Value = 0.62892
$ws.Range("I32").Value = 0.01608
$ws.Range("H33").Value = 0.0887
$ws.Range("I33").Value = 0.02814
$ws.Range("H34").Value = 0.65889
$ws.Range("I34").Value = 0.01267
$ws.Range("H35").Value = 0.08266999999999999
$ws.Range("I35").Value = 0.03898
$ws.Range("H36").Value = 0.65944
$ws.Range("I36").Value = 0.01252
$ws.Range("H37").Value = 0.09213
$ws.Range("I37").Value = 0.04016
$ws.Range("H38").Value = 0.65882
$ws.Range("I38").Value = 0.01285
$ws.Range("H39").Value = 0.08266999999999999
$ws.Range("I39").Value = 0.03898
$ws.Range("H40").Value = 0.6599
$ws.Range("I40").Value = 0.01246
$ws.Range("H41").Value = 0.09213
$ws.Range("I41").Value = 0.04053
$ws.Range("H66").Value = 0.62449
$ws.Range("I66").Value = 0.02246
$ws.Range("H67").Value = 0.10054
$ws.Range("I67").Value = 0.02619
$ws.Range("H68").Value = 0.60481
$ws.Range("I68").Value = 0.0191
$ws.Range("H69").Value = 0.04283
$ws.Range("I69").Value = 0.02219
$ws.Range("H70").Value = 0.62449
$ws.Range("I70").Value = 0.02218
$ws.Range("H71").Value = 0.10122
$ws.Range("I71").Value = 0.02632
$ws.Range("H72").Value = 0.60616
$ws.Range("I72").Value = 0.02466
$ws.Range("H73").Value = 0.05498
$ws.Range("I73").Value = 0.02543
$ws.Range("H74").Value = 0.63514
$ws.Range("I74").Value = 0.02095
$ws.Range("H75").Value = 0.03846
$ws.Range("I75").Value = 0.01677
$ws.Range("H76").Value = 0.6393
$ws.Range("I76").Value = 0.01907
$ws.Range("H77").Value = 0.06815
$ws.Range("I77").Value = 0.02442
$ws.Range("H78").Value = 0.6357
$ws.Range("I78").Value = 0.02067
$ws.Range("H79").Value = 0.03846
$ws.Range("I79").Value = 0.01677
$ws.Range("H80").Value = 0.63964
$ws.Range("I80").Value = 0.01871
$ws.Range("H81").Value = 0.06748
$ws.Range("I81").Value = 0.02444
$ws.Range("H106").Value = 0.66533
$ws.Range("I106").Value = 0.02047
$ws.Range("H107").Value = 0.19329
$ws.Range("I107").Value = 0.0388
$ws.Range("H108").Value = 0.63834
$ws.Range("I108").Value = 0.0207
$ws.Range("H109").Value = 0.0813
$ws.Range("I109").Value = 0.02546
$ws.Range("H110").Value = 0.6654
$ws.Range("I110").Value = 0.0204
$ws.Range("H111").Value = 0.19463
$ws.Range("I111").Value = 0.03993
$ws.Range("H112").Value = 0.63579
$ws.Range("I112").Value = 0.01993
$ws.Range("H113").Value = 0.08569
$ws.Range("I113").Value = 0.02577
$ws.Range("H114").Value = 0.66604
$ws.Range("I114").Value = 0.01348
$ws.Range("H115").Value = 0.07897
$ws.Range("I115").Value = 0.03306
$ws.Range("H116").Value = 0.6662400000000001
$ws.Range("I116").Value = 0.01317
$ws.Range("H117").Value = 0.08705
$ws.Range("I117").Value = 0.03454
$ws.Range("H118").Value = 0.66616
$ws.Range("I118").Value = 0.0135
$ws.Range("H119").Value = 0.07897
$ws.Range("I119").Value = 0.03306
$ws.Range("H120").Value = 0.66636
$ws.Range("I120").Value = 0.01324
$ws.Range("H121").Value = 0.08807
$ws.Range("I121").Value = 0.03491
$ws.Range("H146").Value = 0.62788
$ws.Range("I146").Value = 0.02395
$ws.Range("H147").Value = 0.11811
$ws.Range("I147").Value = 0.03392
$ws.Range("H148").Value = 0.6058
$ws.Range("I148").Value = 0.02255
$ws.Range("H149").Value = 0.05231
$ws.Range("I149").Value = 0.02334
$ws.Range("H150").Value = 0.62927
$ws.Range("I150").Value = 0.02466
$ws.Range("H151").Value = 0.12283
$ws.Range("I151").Value = 0.0355
$ws.Range("H152").Value = 0.60938
$ws.Range("I152").Value = 0.02333
$ws.Range("H153").Value = 0.06884
$ws.Range("I153").Value = 0.02754
$ws.Range("H154").Value = 0.63817
$ws.Range("I154").Value = 0.01607
$ws.Range("H155").Value = 0.0439
$ws.Range("I155").Value = 0.02449
$ws.Range("H156").Value = 0.64352
$ws.Range("I156").Value = 0.01768
$ws.Range("H157").Value = 0.0682
$ws.Range("I157").Value = 0.02925
$ws.Range("H158").Value = 0.63824
$ws.Range("I158").Value = 0.01609
$ws.Range("H159").Value = 0.0439
$ws.Range("I159").Value = 0.02449
$ws.Range("H160").Value = 0.64332
$ws.Range("I160").Value = 0.01808
$ws.Range("H161").Value = 0.06854
$ws.Range("I161").Value = 0.03026
